# Regenerate orders with updated distance/size codes.
# Renames (as substrings, applied across every cell on the sheet):
#   D51 -> D55
#   D64 -> D69
#   D80 -> D86
#   S30 -> S31
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("D51", "D55")
$ws.Cells.Replace("D64", "D69")
$ws.Cells.Replace("D80", "D86")
$ws.Cells.Replace("S30", "S31")
